$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the objective text for the Estimand row (row 3) and CMC row (row 5)
$ws.Range("B3").Value = "The objective of this WG is to promote the concept of estimands and establish a harmonized estimand framework for the clinical development of Cell and Gene Therapies (CGT). Specifically, the WG’s aim is to identify and address the unique challenges posed by intercurrent events in CGT development, and formulate appropriate strategies for their handling."
$ws.Range("B5").Value = "To identify key CMC challenges in CGT development and to serve as key opinion leaders in CGT CMC within biopharmaceutical communities "

# Row 3 grows taller to accommodate the wrapped text
$ws.Rows.Item(3).RowHeight = 45

# Update view: scroll back to column A and move selection
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B14").Select()
